# Actualización automática de noticias - 2026-01-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2..14 down to 3..15.
$ws.Rows.Item(2).Insert()

# The insert copies the bold/centered header formatting down into the new
# row; clear it so the new row matches the plain formatting of the other
# data rows (no explicit style), same as in the target workbook.
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the latest news entry
# Force the "fecha" and "titulo" columns to plain text so values like
# "2026-01-17" and "4" are not reinterpreted as a date/number by Excel.
$ws.Range("A2:B2").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2026-01-17"
$ws.Cells.Item(2, 2).Value = "4"
$ws.Cells.Item(2, 3).Value = "Pulzo"
$ws.Cells.Item(2, 4).Value = "Colombia"
$ws.Cells.Item(2, 5).Value = "https://www.pulzo.com/nacion/que-trata-cambio-estructural-educacion-colombia-por-giro-pae-PP4995431"
$ws.Cells.Item(2, 6).Value = "Anuncian cambio estructural para educación en Colombia: nueva resolución toca a 550.000 niños"
